# Codebook update: add new variable rows (bfi_*_mean, exclude_*, iat_*)
# and widen column A. Existing rows 1-2 are unchanged; rows 3-14 of the
# original sheet are re-seated further down the (now alphabetised) table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@('Variable', 'explanation', 'items_used_and_calculation', 'Range', 'Study')
    ,@('age', 'Participant age', '', '', '')
    ,@('bfi_a_mean', '', '', '', '')
    ,@('bfi_c_mean', '', '', '', '')
    ,@('bfi_e_mean', '', '', '', '')
    ,@('bfi_n_mean', '', '', '', '')
    ,@('bfi_o_mean', '', '', '', '')
    ,@('exclude_bfi_completeness', '', '', '', '')
    ,@('exclude_bfi_logical_scores', '', '', '', '')
    ,@('exclude_iat_completeness', '', '', '', '')
    ,@('exclude_iat_performance', '', '', '', '')
    ,@('exclude_participant', '', '', '', '')
    ,@('exclude_participant_without_bfi', $null, $null, $null, $null)
    ,@('exclude_participant_without_iat', $null, $null, $null, $null)
    ,@('exclude_unique_id_is.na', $null, $null, $null, $null)
    ,@('gender', 'Participant gender', '', '', '')
    ,@('iat_D', '', '', '', '')
    ,@('iat_mean1', '', '', '', '')
    ,@('iat_mean2', '', '', '', '')
    ,@('iat_SD', '', '', '', '')
    ,@('unique_id', 'Participant identification number', '', '', '')
    ,@('bfi_e1,...,bfi_e9', 'Different  items form the subscale extroversion', '', '', '')
    ,@('bfi_c1,...,bfi_c9', 'Different items form the subscale conscientiousness', '', '', '')
    ,@('bfi_n1,...,bfi_n9', 'Different items form the subscale neuroticism', '', '', '')
    ,@('bfi_a1,...,bfi_a9', 'Different items form the subscale agreeableness', '', '', '')
    ,@('bfi_o1,...,bfi_o9', 'Different items form the subscale openness', '', '', '')
    ,@('bfi_e1_rev,...,bfi_e9_rev', 'Different reversed items form the subscale extroversion', '', '', '')
    ,@('bfi_c1_rev,...,bfi_c9_rev', 'Different reversed items form the subscale conscientiousness', '', '', '')
    ,@('bfi_n1_rev,...,bfi_n9_rev', 'Different reversed items form the subscale neuroticism', '', '', '')
    ,@('bfi_a1_rev,...,bfi_a9_rev', 'Different reversed items form the subscale agreeableness', '', '', '')
    ,@('bfi_o1_rev,...,bfi_o9_rev', 'Different reversed items form the subscale openness', '', '', '')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $rowVals = $rows[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $val = $rowVals[$j]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $j + 1).Value = $val
        }
    }
}

# Widen column A (stored width target 31.71 chars; Excel quantises to whole
# pixels, so 30.8 is the closest ColumnWidth that lands on 31.71).
$ws.Columns.Item(1).ColumnWidth = 30.8

Write-Output "codebook rows + column width updated"
